$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 94
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = 46000
$ws.Cells.Item($newRow, 2).Value = "22,1528"
$ws.Cells.Item($newRow, 3).Value = "16,0493"
$ws.Cells.Item($newRow, 4).Value = "15,6461"
$ws.Cells.Item($newRow, 5).Value = "15,6461"

# Match the date/number formatting used by the rest of column A.
$ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($lastRow, 1).NumberFormat
